# Generated PowerShell-style Excel COM-interop script
# Applies the "cryptos list" GitHub Actions update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell as TEXT (never let Excel auto-convert a
# numeric-looking string like "1.02" into a real number), while leaving the
# cells number format back at General so no stray style is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

# --- Rows 44 & 45 swap places (OKB <-> InjectiveProtocol) plus value updates ---
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "25.52"
$ws.Range("E44").Value = "  +1.09%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D45" "41.18"
$ws.Range("E45").Value = "  -0.04%  "

# --- Per-row Price (D) / Volume(1h) (E) updates ---
Set-TextValue "D2" "69.408.26"
$ws.Range("E2").Value = "  +2.02%  "

Set-TextValue "D3" "3.343.27"
$ws.Range("E3").Value = "  +3.01%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D5" "192.27"
$ws.Range("E5").Value = "  +4.08%  "

Set-TextValue "D6" "593.01"
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("E7").Value = "  -0.03%  "

Set-TextValue "D8" "0.608"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("E10").Value = "  +1.63%  "

Set-TextValue "D11" "0.424"
$ws.Range("E11").Value = "  +2.09%  "

Set-TextValue "D12" "3.928.82"
$ws.Range("E12").Value = "  +3.11%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  +1.71%  "

Set-TextValue "D15" "69.424.59"
$ws.Range("E15").Value = "  +2.03%  "

$ws.Range("E16").Value = "  +1.02%  "

Set-TextValue "D17" "3.364.78"
$ws.Range("E17").Value = "  +3.99%  "

$ws.Range("E18").Value = "  +0.42%  "

Set-TextValue "D19" "13.73"
$ws.Range("E19").Value = "  +2.02%  "

Set-TextValue "D20" "427.39"
$ws.Range("E20").Value = "  +7.92%  "

Set-TextValue "D21" "7.71"
$ws.Range("E21").Value = "  +1.59%  "

Set-TextValue "D22" "73.40"
$ws.Range("E22").Value = "  +2.89%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("E26").Value = "  +2.44%  "

Set-TextValue "D27" "9.61"
$ws.Range("E27").Value = "  -0.13%  "

Set-TextValue "D28" "1.02"
$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("E29").Value = "  +2.50%  "

$ws.Range("E30").Value = "  +0.52%  "

Set-TextValue "D31" "23.02"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("E32").Value = "  +1.44%  "

Set-TextValue "D33" "7.01"
$ws.Range("E33").Value = "  -0.12%  "

Set-TextValue "D35" "164.65"
$ws.Range("E35").Value = "  +1.72%  "

Set-TextValue "D36" "1.51"
$ws.Range("E36").Value = "  +1.42%  "

$ws.Range("E37").Value = "  +1.15%  "

Set-TextValue "D38" "27.02"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  -0.16%  "

Set-TextValue "D41" "2.746.23"
$ws.Range("E41").Value = "  +5.10%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("E43").Value = "  +1.13%  "

Set-TextValue "D46" "0.0688"
$ws.Range("E46").Value = "  +0.55%  "

Set-TextValue "D47" "343.93"
$ws.Range("E47").Value = "  +2.56%  "

Set-TextValue "D48" "0.0282"
$ws.Range("E48").Value = "  +1.18%  "

Set-TextValue "D49" "32.49"
$ws.Range("E49").Value = "  +4.46%  "

$ws.Range("E50").Value = "  +3.04%  "

Set-TextValue "D51" "6.29"
$ws.Range("E51").Value = "  -0.33%  "

